$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.131035089492798
$ws.Range("B1").Value = 2.959725141525269
$ws.Range("C1").Value = 2.603036403656006
$ws.Range("D1").Value = 2.830551624298096
$ws.Range("E1").Value = 2.888567447662354
